$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4755.515
$ws.Range("I40").Value = 7325.294
$ws.Range("J40").Value = 2025.125
$ws.Range("K40").Value = 7325.294
$ws.Range("L40").Value = 2025.125
$ws.Range("M40").Value = -7150.294
$ws.Range("N40").Value = -2375.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 386.30435
$ws.Range("I101").Value = 368.06668
$ws.Range("J101").Value = 420.5
$ws.Range("K101").Value = 1104.20004
$ws.Range("L101").Value = 1261.5
$ws.Range("M101").Value = 517.7999599999998
$ws.Range("N101").Value = -4505.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1097472.6
$ws.Range("I132").Value = 1179923.4
$ws.Range("K132").Value = 3539770.2
$ws.Range("M132").Value = -3537240.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2212.4707
$ws.Range("I63").Value = 2124.6667
$ws.Range("J63").Value = 2311.25
$ws.Range("K63").Value = 2124.6667
$ws.Range("L63").Value = 2311.25
$ws.Range("M63").Value = -1438.6667
$ws.Range("N63").Value = -3683.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2212.4707
$ws.Range("I66").Value = 2124.6667
$ws.Range("J66").Value = 2311.25
$ws.Range("K66").Value = 10623.3335
$ws.Range("L66").Value = 11556.25
$ws.Range("M66").Value = -7191.333500000001
$ws.Range("N66").Value = -18420.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 69000
$ws.Range("J101").Value = 69000
$ws.Range("L101").Value = 69000
$ws.Range("N101").Value = -75490

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 22900
$ws.Range("J111").Value = 22900
$ws.Range("L111").Value = 22900
$ws.Range("N111").Value = -31080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 21073
$ws.Range("J123").Value = 21073
$ws.Range("L123").Value = 21073
$ws.Range("N123").Value = -30873

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 35000
$ws.Range("J131").Value = 35000
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 57610.75
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 57610.75
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 57610.75
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -59774.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 50780
$ws.Range("J126").Value = 50780
$ws.Range("L126").Value = 50780
$ws.Range("N126").Value = -60660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5034.2144
$ws.Range("I16").Value = 3934.875
$ws.Range("J16").Value = 6500
$ws.Range("K16").Value = 3934.875
$ws.Range("L16").Value = 6500
$ws.Range("M16").Value = -3647.875
$ws.Range("N16").Value = -7074

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1520.9788
$ws.Range("I31").Value = 1045.1936
$ws.Range("J31").Value = 2442.8125
$ws.Range("K31").Value = 1045.1936
$ws.Range("L31").Value = 2442.8125
$ws.Range("M31").Value = -750.1936000000001
$ws.Range("N31").Value = -3032.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1520.9788
$ws.Range("I34").Value = 1045.1936
$ws.Range("J34").Value = 2442.8125
$ws.Range("K34").Value = 1045.1936
$ws.Range("L34").Value = 2442.8125
$ws.Range("M34").Value = -843.1936000000001
$ws.Range("N34").Value = -2846.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 13335730
$ws.Range("I62").Value = 47621520
$ws.Range("J62").Value = 2366.6667
$ws.Range("K62").Value = 47621520
$ws.Range("L62").Value = 2366.6667
$ws.Range("M62").Value = -47620896
$ws.Range("N62").Value = -3614.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 13335730
$ws.Range("I65").Value = 47621520
$ws.Range("J65").Value = 2366.6667
$ws.Range("K65").Value = 238107600
$ws.Range("L65").Value = 11833.3335
$ws.Range("M65").Value = -238104480
$ws.Range("N65").Value = -18073.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 33447.668
$ws.Range("J88").Value = 33447.668
$ws.Range("L88").Value = 33447.668
$ws.Range("N88").Value = -34259.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 33447.668
$ws.Range("J91").Value = 33447.668
$ws.Range("L91").Value = 33447.668
$ws.Range("N91").Value = -36255.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 29718
$ws.Range("J96").Value = 29718
$ws.Range("L96").Value = 29718
$ws.Range("N96").Value = -35210

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2505.25
$ws.Range("I99").Value = 2212.4707
$ws.Range("J99").Value = 3216.2856
$ws.Range("K99").Value = 2212.4707
$ws.Range("L99").Value = 3216.2856
$ws.Range("M99").Value = -714.4706999999999
$ws.Range("N99").Value = -6212.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 65000
$ws.Range("J106").Value = 65000
$ws.Range("L106").Value = 65000
$ws.Range("N106").Value = -67524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 5034.2144
$ws.Range("I113").Value = 3934.875
$ws.Range("J113").Value = 6500
$ws.Range("K113").Value = 3934.875
$ws.Range("L113").Value = 6500
$ws.Range("M113").Value = -1764.875
$ws.Range("N113").Value = -10840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2505.25
$ws.Range("I126").Value = 2212.4707
$ws.Range("J126").Value = 3216.2856
$ws.Range("K126").Value = 6637.4121
$ws.Range("L126").Value = 9648.856800000001
$ws.Range("M126").Value = -4167.4121
$ws.Range("N126").Value = -14588.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 992.8333
$ws.Range("I69").Value = 877.4
$ws.Range("J69").Value = 1570
$ws.Range("K69").Value = 2632.2
$ws.Range("L69").Value = 4710
$ws.Range("M69").Value = -1821.2
$ws.Range("N69").Value = -6332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 992.8333
$ws.Range("I72").Value = 877.4
$ws.Range("J72").Value = 1570
$ws.Range("K72").Value = 7896.599999999999
$ws.Range("L72").Value = 14130
$ws.Range("M72").Value = -3840.599999999999
$ws.Range("N72").Value = -22242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3620.9092
$ws.Range("I138").Value = 3903.75
$ws.Range("K138").Value = 11711.25
$ws.Range("M138").Value = -6571.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 195.18518
$ws.Range("I55").Value = 158.125
$ws.Range("J55").Value = 249.09091
$ws.Range("K55").Value = 158.125
$ws.Range("L55").Value = 249.09091
$ws.Range("M55").Value = 14.875
$ws.Range("N55").Value = -595.09091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 33495
$ws.Range("J116").Value = 33495
$ws.Range("L116").Value = 33495
$ws.Range("N116").Value = -42673

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2599.8
$ws.Range("I136").Value = 1299.5714
$ws.Range("J136").Value = 3737.5
$ws.Range("K136").Value = 3898.7142
$ws.Range("L136").Value = 11212.5
$ws.Range("M136").Value = -1348.7142
$ws.Range("N136").Value = -16312.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4729.7827
$ws.Range("I107").Value = 372.3
$ws.Range("J107").Value = 8081.6924
$ws.Range("K107").Value = 1116.9
$ws.Range("L107").Value = 24245.0772
$ws.Range("M107").Value = 803.0999999999999
$ws.Range("N107").Value = -28085.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 19000
$ws.Range("J130").Value = 19000
$ws.Range("L130").Value = 19000
$ws.Range("N130").Value = -29040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 26666.666
$ws.Range("J133").Value = 26666.666
$ws.Range("L133").Value = 26666.666
$ws.Range("N133").Value = -36786.666
